$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 05:24:23"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 05:24:17"
$wsZhCn.Range("K2").Value = "2016-09-07 05:24:36"

# de-de sheet: Correspond Handoff Datetime (H2) shares the same original
# text as Overview!G2, so it moves to the same new value; and
# Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 05:24:23"
$wsDeDe.Range("K2").Value = "2016-09-07 05:24:44"
